# Apply "custom accuracy" rounding to row 5 data (from 3 decimals down to 2 decimals)
# and remove the now-redundant row 6, shrinking the used range to A1:AH5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: reduce numeric precision to 2 decimal places (custom accuracy)
$ws.Range("C5").Value  = 8.32
$ws.Range("D5").Value  = 0.94
$ws.Range("E5").Value  = 24.76
$ws.Range("F5").Value  = 20.48
$ws.Range("G5").Value  = 9.07
$ws.Range("H5").Value  = 33.03
$ws.Range("I5").Value  = 13.96
$ws.Range("K5").Value  = 9.15
$ws.Range("L5").Value  = 10.04
$ws.Range("M5").Value  = 10.46
$ws.Range("N5").Value  = 2.9
$ws.Range("O5").Value  = 9.02
$ws.Range("P5").Value  = 12.75
$ws.Range("Q5").Value  = 7.71
$ws.Range("S5").Value  = 0.57
$ws.Range("T5").Value  = 129.99
$ws.Range("U5").Value  = 25.15
$ws.Range("V5").Value  = 8.33
$ws.Range("W5").Value  = 16.76
$ws.Range("X5").Value  = 9.02
$ws.Range("Y5").Value  = 1.07
$ws.Range("Z5").Value  = 16.15
$ws.Range("AA5").Value = 7.36
$ws.Range("AB5").Value = 6.6
$ws.Range("AC5").Value = 7.75
$ws.Range("AD5").Value = 10.48
$ws.Range("AE5").Value = 0.56
$ws.Range("AF5").Value = 29.61
$ws.Range("AG5").Value = 4.63
$ws.Range("AH5").Value = 10.41

# Row 6 is removed entirely (data trimmed), which also shrinks the sheet
# dimension from A1:AH6 down to A1:AH5.
$ws.Rows("6:6").Delete()
